# Apply corrected Diebold-Mariano statistics (DM_Stat / P_Value) for each
# comparison row. Only columns C (DM_Stat) and D (P_Value) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ C = 1.050264827978696;  D = 0.3010044991370782 }
    3  = @{ C = 1.194885992271794;  D = 0.2404003166471858 }
    4  = @{ C = 0.9842762547185959; D = 0.3319325096132584 }
    5  = @{ C = 0.6720168418813248; D = 0.5061141733098915 }
    6  = @{ C = -0.1817074654141688; D = 0.8568912969935161 }
    7  = @{ C = -0.2874160849894225; D = 0.7755383101645201 }
    8  = @{ C = -0.3370678349915029; D = 0.738137713428662 }
    9  = @{ C = -0.0138307112247822; D = 0.9890458813589216 }
    10 = @{ C = -0.07025206961457434; D = 0.9444046199335046 }
    11 = @{ C = -0.1198031897139542; D = 0.9053441842806356 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("D$row").Value = $values[$row].D
}
